$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix provider/contractor names: comma separators mistakenly used -> periods ---
$ws.Cells.Item(50,5).Value = "ALBIZZATTI. PABLO MARTIN Y FULINI. SERGIO RUBEN"
$ws.Cells.Item(67,5).Value = "MARSICO GUILLERMO MIGUEL. MARSICO JUAN EDUARDO"
$ws.Cells.Item(104,5).Value = "MARSICO GUILLERMO MIGUEL. MARSICO JUAN EDUARDO"
$ws.Cells.Item(185,5).Value = "RICCOTTI. MARIANA EDITH"
$ws.Cells.Item(205,5).Value = "ALBIZZATTI. PABLO MARTIN Y FULINI. SERGIO RUBEN"

# --- Fix "Importe" (amount) column H: convert from es-AR "1.234,56" text to "1234.56" ---
$amountRange = $ws.Range("H2:H250")
$amountRange.NumberFormat = "@"

$ws.Cells.Item(2,8).Value = "3140.00"
$ws.Cells.Item(3,8).Value = "1500.00"
$ws.Cells.Item(4,8).Value = "1705.00"
$ws.Cells.Item(5,8).Value = "300000.00"
$ws.Cells.Item(6,8).Value = "25680.00"
$ws.Cells.Item(7,8).Value = "9376.60"
$ws.Cells.Item(8,8).Value = "1790.00"
$ws.Cells.Item(9,8).Value = "5980.00"
$ws.Cells.Item(10,8).Value = "7020.00"
$ws.Cells.Item(11,8).Value = "93000.00"
$ws.Cells.Item(12,8).Value = "583.60"
$ws.Cells.Item(13,8).Value = "3680.00"
$ws.Cells.Item(14,8).Value = "110.00"
$ws.Cells.Item(15,8).Value = "1716.00"
$ws.Cells.Item(16,8).Value = "1100.00"
$ws.Cells.Item(17,8).Value = "110117.92"
$ws.Cells.Item(18,8).Value = "18000.00"
$ws.Cells.Item(19,8).Value = "7712.50"
$ws.Cells.Item(20,8).Value = "12480.00"
$ws.Cells.Item(21,8).Value = "16689.79"
$ws.Cells.Item(22,8).Value = "459.00"
$ws.Cells.Item(23,8).Value = "70547.96"
$ws.Cells.Item(24,8).Value = "1098.00"
$ws.Cells.Item(25,8).Value = "18489.00"
$ws.Cells.Item(26,8).Value = "46828.85"
$ws.Cells.Item(27,8).Value = "50250.00"
$ws.Cells.Item(28,8).Value = "9000.00"
$ws.Cells.Item(29,8).Value = "6800.00"
$ws.Cells.Item(30,8).Value = "3500.00"
$ws.Cells.Item(31,8).Value = "7200.00"
$ws.Cells.Item(32,8).Value = "99.20"
$ws.Cells.Item(33,8).Value = "500.00"
$ws.Cells.Item(34,8).Value = "45.00"
$ws.Cells.Item(35,8).Value = "70.00"
$ws.Cells.Item(36,8).Value = "140181.94"
$ws.Cells.Item(37,8).Value = "3262.34"
$ws.Cells.Item(38,8).Value = "101858.96"
$ws.Cells.Item(39,8).Value = "965.58"
$ws.Cells.Item(40,8).Value = "2585.80"
$ws.Cells.Item(41,8).Value = "8980.00"
$ws.Cells.Item(42,8).Value = "5813.16"
$ws.Cells.Item(43,8).Value = "300.00"
$ws.Cells.Item(44,8).Value = "390.07"
$ws.Cells.Item(45,8).Value = "32888.32"
$ws.Cells.Item(46,8).Value = "3365.00"
$ws.Cells.Item(47,8).Value = "1400.00"
$ws.Cells.Item(48,8).Value = "66294.74"
$ws.Cells.Item(49,8).Value = "249.03"
$ws.Cells.Item(50,8).Value = "1740.00"
$ws.Cells.Item(51,8).Value = "86952.48"
$ws.Cells.Item(52,8).Value = "187.00"
$ws.Cells.Item(53,8).Value = "8300.00"
$ws.Cells.Item(54,8).Value = "300.00"
$ws.Cells.Item(55,8).Value = "2218.90"
$ws.Cells.Item(56,8).Value = "42.02"
$ws.Cells.Item(57,8).Value = "2096.00"
$ws.Cells.Item(58,8).Value = "3917.80"
$ws.Cells.Item(59,8).Value = "15507.10"
$ws.Cells.Item(60,8).Value = "300.00"
$ws.Cells.Item(61,8).Value = "8422.00"
$ws.Cells.Item(62,8).Value = "17681.84"
$ws.Cells.Item(63,8).Value = "3044.08"
$ws.Cells.Item(64,8).Value = "240.00"
$ws.Cells.Item(65,8).Value = "2900.00"
$ws.Cells.Item(66,8).Value = "201.13"
$ws.Cells.Item(67,8).Value = "868.00"
$ws.Cells.Item(68,8).Value = "9200.00"
$ws.Cells.Item(69,8).Value = "233497.02"
$ws.Cells.Item(70,8).Value = "1240.00"
$ws.Cells.Item(71,8).Value = "7188.00"
$ws.Cells.Item(72,8).Value = "968.00"
$ws.Cells.Item(73,8).Value = "420.00"
$ws.Cells.Item(74,8).Value = "65.00"
$ws.Cells.Item(75,8).Value = "4225.16"
$ws.Cells.Item(76,8).Value = "12538.33"
$ws.Cells.Item(77,8).Value = "9546.00"
$ws.Cells.Item(78,8).Value = "16162.90"
$ws.Cells.Item(79,8).Value = "9880.00"
$ws.Cells.Item(80,8).Value = "1535.00"
$ws.Cells.Item(81,8).Value = "4000.00"
$ws.Cells.Item(82,8).Value = "14895.00"
$ws.Cells.Item(83,8).Value = "17500.00"
$ws.Cells.Item(84,8).Value = "25000.00"
$ws.Cells.Item(85,8).Value = "278850.00"
$ws.Cells.Item(86,8).Value = "11100.00"
$ws.Cells.Item(87,8).Value = "18000.00"
$ws.Cells.Item(88,8).Value = "3034.72"
$ws.Cells.Item(89,8).Value = "518.00"
$ws.Cells.Item(90,8).Value = "763.00"
$ws.Cells.Item(91,8).Value = "1040.00"
$ws.Cells.Item(92,8).Value = "6812.60"
$ws.Cells.Item(93,8).Value = "2160.00"
$ws.Cells.Item(94,8).Value = "16878.00"
$ws.Cells.Item(95,8).Value = "2600.00"
$ws.Cells.Item(96,8).Value = "480.00"
$ws.Cells.Item(97,8).Value = "3300.00"
$ws.Cells.Item(98,8).Value = "50000.00"
$ws.Cells.Item(99,8).Value = "1000.00"
$ws.Cells.Item(100,8).Value = "4800.00"
$ws.Cells.Item(101,8).Value = "700.00"
$ws.Cells.Item(102,8).Value = "95.00"
$ws.Cells.Item(103,8).Value = "6340.00"
$ws.Cells.Item(104,8).Value = "2750.00"
$ws.Cells.Item(105,8).Value = "51360.00"
$ws.Cells.Item(106,8).Value = "3000.00"
$ws.Cells.Item(107,8).Value = "4221.00"
$ws.Cells.Item(108,8).Value = "13.78"
$ws.Cells.Item(109,8).Value = "132.00"
$ws.Cells.Item(110,8).Value = "63600.00"
$ws.Cells.Item(111,8).Value = "275.59"
$ws.Cells.Item(112,8).Value = "138.00"
$ws.Cells.Item(113,8).Value = "27186.42"
$ws.Cells.Item(114,8).Value = "707.00"
$ws.Cells.Item(115,8).Value = "81.00"
$ws.Cells.Item(116,8).Value = "1380.00"
$ws.Cells.Item(117,8).Value = "1855.10"
$ws.Cells.Item(118,8).Value = "3430.00"
$ws.Cells.Item(119,8).Value = "1631.00"
$ws.Cells.Item(120,8).Value = "7599.00"
$ws.Cells.Item(121,8).Value = "4753.00"
$ws.Cells.Item(122,8).Value = "450.00"
$ws.Cells.Item(123,8).Value = "2395.00"
$ws.Cells.Item(124,8).Value = "7252.00"
$ws.Cells.Item(125,8).Value = "8047.00"
$ws.Cells.Item(126,8).Value = "81508.01"
$ws.Cells.Item(127,8).Value = "2548.70"
$ws.Cells.Item(128,8).Value = "12289.20"
$ws.Cells.Item(129,8).Value = "5594.50"
$ws.Cells.Item(130,8).Value = "2840.00"
$ws.Cells.Item(131,8).Value = "2571.50"
$ws.Cells.Item(132,8).Value = "329.50"
$ws.Cells.Item(133,8).Value = "740.00"
$ws.Cells.Item(134,8).Value = "41400.00"
$ws.Cells.Item(135,8).Value = "5316.00"
$ws.Cells.Item(136,8).Value = "278.00"
$ws.Cells.Item(137,8).Value = "1760.00"
$ws.Cells.Item(138,8).Value = "1136.00"
$ws.Cells.Item(139,8).Value = "180.00"
$ws.Cells.Item(140,8).Value = "4700.00"
$ws.Cells.Item(141,8).Value = "40000.00"
$ws.Cells.Item(142,8).Value = "439710.00"
$ws.Cells.Item(143,8).Value = "160.00"
$ws.Cells.Item(144,8).Value = "2995.00"
$ws.Cells.Item(145,8).Value = "7072.92"
$ws.Cells.Item(146,8).Value = "7160.00"
$ws.Cells.Item(147,8).Value = "740.00"
$ws.Cells.Item(148,8).Value = "3500.00"
$ws.Cells.Item(149,8).Value = "106660.00"
$ws.Cells.Item(150,8).Value = "4800.00"
$ws.Cells.Item(151,8).Value = "48000.00"
$ws.Cells.Item(152,8).Value = "10000.00"
$ws.Cells.Item(153,8).Value = "19000.00"
$ws.Cells.Item(154,8).Value = "18400.00"
$ws.Cells.Item(155,8).Value = "80500.00"
$ws.Cells.Item(156,8).Value = "26300.00"
$ws.Cells.Item(157,8).Value = "2747.96"
$ws.Cells.Item(158,8).Value = "1678.00"
$ws.Cells.Item(159,8).Value = "461.13"
$ws.Cells.Item(160,8).Value = "3400.00"
$ws.Cells.Item(161,8).Value = "477600.00"
$ws.Cells.Item(162,8).Value = "50000.00"
$ws.Cells.Item(163,8).Value = "12500.00"
$ws.Cells.Item(164,8).Value = "14000.00"
$ws.Cells.Item(165,8).Value = "18000.00"
$ws.Cells.Item(166,8).Value = "8000.00"
$ws.Cells.Item(167,8).Value = "6500.00"
$ws.Cells.Item(168,8).Value = "5000.00"
$ws.Cells.Item(169,8).Value = "8508.50"
$ws.Cells.Item(170,8).Value = "12000.00"
$ws.Cells.Item(171,8).Value = "7000.00"
$ws.Cells.Item(172,8).Value = "5000.00"
$ws.Cells.Item(173,8).Value = "6000.00"
$ws.Cells.Item(174,8).Value = "6000.00"
$ws.Cells.Item(175,8).Value = "5000.00"
$ws.Cells.Item(176,8).Value = "5000.00"
$ws.Cells.Item(177,8).Value = "12000.00"
$ws.Cells.Item(178,8).Value = "7000.00"
$ws.Cells.Item(179,8).Value = "12500.00"
$ws.Cells.Item(180,8).Value = "6000.00"
$ws.Cells.Item(181,8).Value = "8000.00"
$ws.Cells.Item(182,8).Value = "6500.00"
$ws.Cells.Item(183,8).Value = "20000.00"
$ws.Cells.Item(184,8).Value = "45000.00"
$ws.Cells.Item(185,8).Value = "12000.00"
$ws.Cells.Item(186,8).Value = "6500.00"
$ws.Cells.Item(187,8).Value = "6000.00"
$ws.Cells.Item(188,8).Value = "46161.50"
$ws.Cells.Item(189,8).Value = "3500.00"
$ws.Cells.Item(190,8).Value = "38608.00"
$ws.Cells.Item(191,8).Value = "6914.48"
$ws.Cells.Item(192,8).Value = "26900.00"
$ws.Cells.Item(193,8).Value = "4200.00"
$ws.Cells.Item(194,8).Value = "90800.00"
$ws.Cells.Item(195,8).Value = "74450.00"
$ws.Cells.Item(196,8).Value = "4200.00"
$ws.Cells.Item(197,8).Value = "3200.00"
$ws.Cells.Item(198,8).Value = "1050.00"
$ws.Cells.Item(199,8).Value = "9450.00"
$ws.Cells.Item(200,8).Value = "9564.00"
$ws.Cells.Item(201,8).Value = "7705.50"
$ws.Cells.Item(202,8).Value = "110.00"
$ws.Cells.Item(203,8).Value = "8350.00"
$ws.Cells.Item(204,8).Value = "1152.60"
$ws.Cells.Item(205,8).Value = "1577.00"
$ws.Cells.Item(206,8).Value = "29585.00"
$ws.Cells.Item(207,8).Value = "60914.00"
$ws.Cells.Item(208,8).Value = "14512.00"
$ws.Cells.Item(209,8).Value = "1324.84"
$ws.Cells.Item(210,8).Value = "5200.00"
$ws.Cells.Item(211,8).Value = "1359.48"
$ws.Cells.Item(212,8).Value = "12575.92"
$ws.Cells.Item(213,8).Value = "5666.03"
$ws.Cells.Item(214,8).Value = "10925.36"
$ws.Cells.Item(215,8).Value = "3455.76"
$ws.Cells.Item(216,8).Value = "27707.62"
$ws.Cells.Item(217,8).Value = "450.00"
$ws.Cells.Item(218,8).Value = "2069.10"
$ws.Cells.Item(219,8).Value = "3000.00"
$ws.Cells.Item(220,8).Value = "165.00"
$ws.Cells.Item(221,8).Value = "6616.16"
$ws.Cells.Item(222,8).Value = "3369.95"
$ws.Cells.Item(223,8).Value = "90.00"
$ws.Cells.Item(224,8).Value = "57000.00"
$ws.Cells.Item(225,8).Value = "2500.00"
$ws.Cells.Item(226,8).Value = "8862.00"
$ws.Cells.Item(227,8).Value = "844173.37"
$ws.Cells.Item(228,8).Value = "5198.16"
$ws.Cells.Item(229,8).Value = "6719055.13"
$ws.Cells.Item(230,8).Value = "7200.00"
$ws.Cells.Item(231,8).Value = "184321.00"
$ws.Cells.Item(232,8).Value = "1132413.37"
$ws.Cells.Item(233,8).Value = "9000.00"
$ws.Cells.Item(234,8).Value = "3450.00"
$ws.Cells.Item(235,8).Value = "6000.00"
$ws.Cells.Item(236,8).Value = "4750.00"
$ws.Cells.Item(237,8).Value = "51900.79"
$ws.Cells.Item(238,8).Value = "29100.00"
$ws.Cells.Item(239,8).Value = "51000.00"
$ws.Cells.Item(240,8).Value = "47000.00"
$ws.Cells.Item(241,8).Value = "4000.00"
$ws.Cells.Item(242,8).Value = "109162.16"
$ws.Cells.Item(243,8).Value = "752.39"
$ws.Cells.Item(244,8).Value = "800.00"
$ws.Cells.Item(245,8).Value = "9000.00"
$ws.Cells.Item(246,8).Value = "8330.00"
$ws.Cells.Item(247,8).Value = "49200.00"
$ws.Cells.Item(248,8).Value = "31812.95"
$ws.Cells.Item(249,8).Value = "47000.00"
$ws.Cells.Item(250,8).Value = "28121.00"

$amountRange.Style = "Normal"
